{"js": "// Apply the \"DB linked, Server Started, First route\" edit.\n//\n// 1. Merge the two adjacent (identically-formatted) runs that make up\n//    \" npm i cors \u2013save\" into a single run (Word silently coalesces\n//    runs like this when the text is retyped / re-saved).\n// 2. Drop the leftover \"_GoBack\" bookmark (Word regenerates / discards\n//    this automatically; it marks the last edit location).\n// 3. Add the new paragraph documenting how to create the MongoDB\n//    database, right after the \"Generar Conexi\u00f3n a la base de datos de\n//    Mongoose\" heading (i.e. immediately before the trailing empty\n//    paragraph that closes the document body).\n\nconst body = context.document.body;\n\n// --- 1. Merge the \"npm i cors \u2013save\" runs -------------------------------\nconst corsResults = body.search(\" npm i cors \u2013save\", { matchCase: true });\ncorsResults.load(\"items\");\nawait context.sync();\n\nif (corsResults.items.length > 0) {\n  // Re-inserting the identical text over the matched range causes the\n  // host to rebuild it as a single run instead of the two runs\n  // (\" \" + \"npm i cors \u2013save\") the source document had.\n  corsResults.items[0].insertText(\" npm i cors \u2013save\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Remove the stale \"_GoBack\" bookmark ------------------------------\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 3. Insert the new \"Crear Base de datos en MongoDB...\" paragraph ----\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst heading = items.find(\n  (p) => p.text === \"Generar Conexi\u00f3n a la base de datos de Mongoose\"\n);\n\n// The document body ends with: \"...Generar Conexi\u00f3n...\" heading, then one\n// trailing empty paragraph before the section break. Insert the new\n// paragraph right before that trailing paragraph (i.e. right after the\n// heading) so it naturally picks up the same (non-bold, 12pt, text1) run\n// formatting instead of the bold 16pt heading formatting.\nconst insertionPoint = heading\n  ? heading.getNextOrNullObject()\n  : null;\ninsertionPoint && insertionPoint.load(\"isNullObject\");\nawait context.sync();\n\nconst target =\n  insertionPoint && !insertionPoint.isNullObject\n    ? insertionPoint\n    : items[items.length - 1];\n\ntarget.insertParagraph(\n  \"Crear Base de datos en MongoDB, usar cmd o MongoCompass, creamos la base de datos mi_blog con una colecci\u00f3n llamada articles.\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n", "ps1": "# Apply the \"DB linked, Server Started, First route\" edit.\n#\n# 1. Merge the two adjacent (identically-formatted) runs that make up\n#    \" npm i cors \u2013save\" into a single run (Word silently coalesces\n#    runs like this when the text is retyped / re-saved).\n# 2. Drop the leftover \"_GoBack\" bookmark (Word regenerates / discards\n#    this automatically; it marks the last edit location).\n# 3. Add the new paragraph documenting how to create the MongoDB\n#    database, right after the \"Generar Conexi\u00f3n a la base de datos de\n#    Mongoose\" heading (i.e. immediately before the trailing empty\n#    paragraph that closes the document body).\n\n$d = $word.ActiveDocument\n\n# --- 1. Merge the \"npm i cors \u2013save\" runs --------------------------------\n$findRange = $d.Content\n$findRange.Find.Text = \" npm i cors \u2013save\"\n$findRange.Find.Replacement.Text = \" npm i cors \u2013save\"\n$null = $findRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# --- 2. Remove the stale \"_GoBack\" bookmark ------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3. Insert the new \"Crear Base de datos en MongoDB...\" paragraph ----\n# The document body ends with: \"...Generar Conexi\u00f3n...\" heading, then one\n# trailing empty paragraph before the section break. Insert the new\n# paragraph right before that trailing paragraph (i.e. right after the\n# heading) so it naturally picks up the same (non-bold, 12pt, text1) run\n# formatting instead of the bold 16pt heading formatting.\n$headingPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Generar Conexi\u00f3n a la base de datos de Mongoose\") {\n        $headingPara = $p\n        break\n    }\n}\n\nif ($headingPara -ne $null) {\n    $target = $headingPara.Next().Range\n} else {\n    # Fallback: couldn't locate the heading, just land on the last paragraph.\n    $target = $d.Paragraphs.Last.Range\n}\n\n$target.Collapse(1)\n$target.InsertBefore(\"Crear Base de datos en MongoDB, usar cmd o MongoCompass, creamos la base de datos mi_blog con una colecci\u00f3n llamada articles.\" + [char]13)\n"}
